$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "October 15 2023" meeting-minutes entry in row 10 -----------
# Mirror the formatting already used by the equivalent cells in row 9
# (Date -> style like A9, Attendance -> style like B9, Started -> style like C9)
# and fill in the new meeting's info, matching the existing table layout:
# A=Date, B=Attendance, C=Started, D=Ended, E=What was Done, F=What was Confirmed

$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A10").Value = "October 15 2023"

$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C10").Value = "7:00PM"

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B10").Value = "Sedat, Madison, Joseph, David"

$excel.CutCopyMode = 0

# --- Widen column E slightly (as it was resized in the original edit) -----
$ws.Range("E:E").ColumnWidth = 68.916666666666671

# --- Move the active selection to the newly-edited row ---------------------
$ws.Range("B10").Select()
